# Apply cryptos list update (Tue May 16 11:22:04 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.156.45'
$ws.Range("E2").Value = '  -1.99%  '

$ws.Range("D3").Value = '1.821.26'
$ws.Range("E3").Value = '  -1.49%  '

$c = $ws.Range("D4")
$c.Value = "'1.005"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.88%  '

$c = $ws.Range("D5")
$c.Value = "'312.13"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.34%  '

$c = $ws.Range("D6")
$c.Value = "'1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.78%  '

$c = $ws.Range("D7")
$c.Value = "'0.4226"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.12%  '

$c = $ws.Range("D8")
$c.Value = "'0.3678"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.74%  '

$c = $ws.Range("D9")
$c.Value = "'0.07234"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.64%  '

$c = $ws.Range("D10")
$c.Value = "'0.8540"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.99%  '

$c = $ws.Range("D11")
$c.Value = "'20.95"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.06%  '

$ws.Range("D12").Value = '1.827.23'
$ws.Range("E12").Value = '  -0.93%  '

$c = $ws.Range("D13")
$c.Value = "'6.692"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.60%  '

$c = $ws.Range("D14")
$c.Value = "'0.07071"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.66%  '

$c = $ws.Range("D15")
$c.Value = "'5.289"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.02%  '

$c = $ws.Range("D16")
$c.Value = "'89.61"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.81%  '

$ws.Range("E17").Value = '  -0.87%  '

$c = $ws.Range("D18")
$c.Value = "'0.000008838"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.74%  '

$ws.Range("E19").Value = '  -0.74%  '

$c = $ws.Range("D20")
$c.Value = "'15.00"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.06%  '

$ws.Range("D21").Value = '27.228.69'
$ws.Range("E21").Value = '  -1.72%  '

$c = $ws.Range("D22")
$c.Value = "'5.110"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.64%  '

$c = $ws.Range("D23")
$c.Value = "'10.85"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.82%  '

$ws.Range("D24").Value = '2.054.62'
$ws.Range("E24").Value = '  -0.89%  '

$c = $ws.Range("D25")
$c.Value = "'1.978"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '

$c = $ws.Range("D26")
$c.Value = "'152.31"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.15%  '

$c = $ws.Range("D27")
$c.Value = "'2.197"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.83%  '

$c = $ws.Range("D28")
$c.Value = "'18.35"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.42%  '

$c = $ws.Range("D29")
$c.Value = "'5.226"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.12%  '

$c = $ws.Range("D30")
$c.Value = "'116.06"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.72%  '

$c = $ws.Range("D31")
$c.Value = "'0.08831"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.21%  '

$c = $ws.Range("D32")
$c.Value = "'1.188"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.36%  '

$c = $ws.Range("D33")
$c.Value = "'0.7467"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.12%  '

$c = $ws.Range("D34")
$c.Value = "'2.945"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.94%  '

$ws.Range("E35").Value = '  -2.84%  '

$c = $ws.Range("D36")
$c.Value = "'1.003"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.73%  '

$c = $ws.Range("D37")
$c.Value = "'1.107"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.01%  '

$c = $ws.Range("D38")
$c.Value = "'0.01963"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.52%  '

$c = $ws.Range("D39")
$c.Value = "'0.05232"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.88%  '

$c = $ws.Range("D40")
$c.Value = "'7.265"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.18%  '

$c = $ws.Range("D41")
$c.Value = "'2.861"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '

$c = $ws.Range("D42")
$c.Value = "'0.1694"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.89%  '

$c = $ws.Range("D43")
$c.Value = "'0.5024"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.63%  '

$c = $ws.Range("D44")
$c.Value = "'8.636"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.15%  '

$c = $ws.Range("D45")
$c.Value = "'10.57"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.31%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D46")
$c.Value = "'106.27"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.99%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D47")
$c.Value = "'0.4733"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.01%  '

$c = $ws.Range("D48")
$c.Value = "'1.003"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.81%  '

$c = $ws.Range("D49")
$c.Value = "'0.06386"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.94%  '

$c = $ws.Range("D50")
$c.Value = "'1.658"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.43%  '

$c = $ws.Range("D51")
$c.Value = "'1.873"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.40%  '
